# Append two new "news" rows (39 and 40) to the Historico sheet,
# matching the newsbot state update described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 39 --------------------------------------------------------------
$v39A = @'
05/01/2026 04:06:31
'@
$ws.Range("A39").Value = $v39A

$v39B = @'
05/01 04:01
'@
$ws.Range("B39").Value = $v39B

$v39C = @'
g1 > Economia
'@
$ws.Range("C39").Value = $v39C

$v39D = @'
Inflação dentro da meta? Veja os preços que mais caíram e os que mais subiram em 2025
'@
$ws.Range("D39").Value = $v39D

$v39E = @'
https://g1.globo.com/economia/noticia/2026/01/05/veja-os-precos-que-mais-cairam-e-os-que-mais-subiram-em-2025.ghtml
'@
$ws.Range("E39").Value = $v39E

$v39F = @'
orçamento
'@
$ws.Range("F39").Value = $v39F

$v39G = @'
nsino fundamental;
Empregado doméstico; e 
Condomínio.
“Juntos, eles representam 15,8% do &lt;b&gt;orçamento&lt;/b&gt; doméstico e registraram uma inflação média de 6,2% entre janeiro e novembro de 2025, acim
'@
$ws.Range("G39").Value = $v39G

# --- Row 40 --------------------------------------------------------------
$v40A = @'
05/01/2026 04:06:32
'@
$ws.Range("A40").Value = $v40A

$v40B = @'
05/01 04:01
'@
$ws.Range("B40").Value = $v40B

$v40C = @'
g1 > Economia
'@
$ws.Range("C40").Value = $v40C

$v40D = @'
ES vai receber R$ 106 bi em investimentos até 2029; veja obras e valores por cidade
'@
$ws.Range("D40").Value = $v40D

$v40E = @'
https://g1.globo.com/es/espirito-santo/noticia/2026/01/05/es-vai-receber-r-106-bi-em-investimentos-ate-2029-veja-obras-e-valores-por-cidade.ghtml
'@
$ws.Range("E40").Value = $v40E

$v40F = @'
pix
'@
$ws.Range("F40").Value = $v40F

$v40G = @'
_photos/bs/2025/P/J/c2IHCASrAuLKmuUfOnkg/thumb-videos-g1-21-.jpg" /&gt;&lt;br /&gt;     Economia ca&lt;b&gt;pix&lt;/b&gt;aba cresce 2,2% em nove meses puxada pela agropecuária e indústria
O Espírito Santo vai re
'@
$ws.Range("G40").Value = $v40G

